# Rows 30-34 hold 5 sightings that get re-ordered/re-keyed (Id, Ost/Nord
# coordinates, times, species, etc. rotate down by one row, with the last
# row's original content wrapping around to the first row).
#
# Net effect verified against the target diff:
#   new row 30 = old row 31
#   new row 31 = old row 32
#   new row 32 = old row 33
#   new row 33 = old row 34
#   new row 34 = old row 30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30 (becomes old row 31's data: Tretåig hackspett / gran, 09:47) ---
$ws.Range("A30").Value = 130979946
$ws.Range("B30").Value = 57884
$ws.Range("E30").Value = 100109
$ws.Range("F30").Value = "Tretåig hackspett"
$ws.Range("G30").Value = "Picoides tridactylus"
$ws.Range("H30").Value = "(Linnaeus, 1758)"
$ws.Range("K30").Value = ""
$ws.Range("L30").Value = ""
$ws.Range("M30").Value = "färska spår"
$ws.Range("N30").Value = ""
$ws.Range("Q30").Value = 590605
$ws.Range("R30").Value = 6963364
$ws.Range("Z30").Value = "09:47"
$ws.Range("AB30").Value = "09:47"
$ws.Range("AC30").Value = "färska ringhack på gran"

# --- Row 31 (becomes old row 32's data: Tretåig hackspett / gran, 13:16) ---
$ws.Range("A31").Value = 130979899
$ws.Range("Q31").Value = 590850
$ws.Range("R31").Value = 6963133
$ws.Range("Z31").Value = "13:16"
$ws.Range("AB31").Value = "13:16"

# --- Row 32 (becomes old row 33's data: Tretåig hackspett / tall, 11:44) ---
$ws.Range("A32").Value = 130979914
$ws.Range("Q32").Value = 591126
$ws.Range("R32").Value = 6963169
$ws.Range("Z32").Value = "11:44"
$ws.Range("AB32").Value = "11:44"
$ws.Range("AC32").Value = "färska ringhack på tall"

# --- Row 33 (becomes old row 34's data: Lunglav, 13:24, no K/L/M/N/AC) ---
$ws.Range("A33").Value = 130979897
$ws.Range("B33").Value = 80348
$ws.Range("E33").Value = 6458
$ws.Range("F33").Value = "Lunglav"
$ws.Range("G33").Value = "Lobaria pulmonaria"
$ws.Range("H33").Value = "(L.) Hoffm."
$ws.Range("K33").ClearContents()
$ws.Range("L33").ClearContents()
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()
$ws.Range("Q33").Value = 590726
$ws.Range("R33").Value = 6963153
$ws.Range("Z33").Value = "13:24"
$ws.Range("AB33").Value = "13:24"
$ws.Range("AC33").ClearContents()

# --- Row 34 (becomes old row 30's data: Ullticka, 09:45) ---
$ws.Range("A34").Value = 130979947
$ws.Range("B34").Value = 91808
$ws.Range("E34").Value = 1202
$ws.Range("F34").Value = "Ullticka"
$ws.Range("G34").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H34").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q34").Value = 590591
$ws.Range("R34").Value = 6963354
$ws.Range("Z34").Value = "09:45"
$ws.Range("AB34").Value = "09:45"
